$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 25.73000000000058
$ws.Range("H2").Value = [double]"5.657187386625001e-16"
$ws.Range("K2").Value = 48.04868932699065
$ws.Range("L2").Value = "[39.24213368426864, 56.85524496971266]"
$ws.Range("O2").Value = 1.66671081721981
$ws.Range("P2").Value = "[1.46544762419704, 1.8679740102425804]"
$ws.Range("S2").Value = 65.41415439512311
$ws.Range("T2").Value = "[59.85210491132517, 70.97620387892104]"
$ws.Range("W2").Value = 18.90472472472515
$ws.Range("X2").Value = 18.08054054054095
$ws.Range("Y2").Value = 19.72890890890936

# Row 3 updates
$ws.Range("E3").Value = 24.5800000000004
$ws.Range("G3").Value = [double]"1.46960221769632e-12"
$ws.Range("H3").Value = [double]"5.367860351224389e-12"
$ws.Range("I3").Value = "'"
$ws.Range("I3").ClearFormats()
$ws.Range("K3").Value = 51.84731901909259
$ws.Range("L3").Value = "[35.50179658660899, 68.1928414515762]"
$ws.Range("M3").Value = [double]"4.125253916242855e-09"
$ws.Range("N3").Value = [double]"4.125253916242855e-09"
$ws.Range("O3").Value = -0.5031579825569237
$ws.Range("P3").Value = "[-0.8427896207828471, -0.16352634433100022]"
$ws.Range("Q3").Value = 0.003969598243227956
$ws.Range("R3").Value = 0.003969598243227956
$ws.Range("S3").Value = 62.73509323844965
$ws.Range("T3").Value = "[53.760133263172065, 71.71005321372724]"
$ws.Range("W3").Value = 1.968368368368399
$ws.Range("X3").Value = 0.6397197197197289
$ws.Range("Y3").Value = 3.297017017017069
